$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
try {
    $win.TopLeftCell = $ws.Range("A128")
    Write-Host "set TopLeftCell via window ok"
} catch {
    Write-Host "window.TopLeftCell failed:" $_
}
try {
    $v = $win.ScrollRow()
    Write-Host "ScrollRow:" $v
} catch {
    Write-Host "ScrollRow read failed:" $_
}
try {
    $panes = $win.Panes
    Write-Host "panes:" $panes
    Write-Host "panes count:" $panes.Count()
} catch {
    Write-Host "Panes failed:" $_
}
